$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.382, 0.537, -0.121, 140, -34, 36),
    @(-0.508, -0.249, 0.233, 141, -10, 127),
    @(0.095, -0.262, -0.276, -148, -79, -67),
    @(-0.596, -0.473, -0.035, 111, -27, -36),
    @(-0.005, -0.152, 0.241, -128, 7, 64),
    @(-0.151, 0.021, -0.358, -68, -17, -75),
    @(0.008999999999999999, 0.29, -0.346, -105, -73, 118),
    @(0.214, 0.466, -0.063, -123, -52, 55),
    @(0.218, -0.152, 0.582, -143, -79, -125),
    @(-0.524, -0.475, -0.028, -26, -44, -167)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}
